$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ExtremeFlows")

# Update the "Extreme Low Flow Method" text for each row, appending a
# trailing period to the existing entries and replacing the A2 method
# description with the new wording. Row positions stay the same; only
# the text (and, for row 3, the Minimum/Maximum values) change.
$ws.Range("A2").Value = "A1. Lowest 10-year average flows:in Reclamation's post-:2026 ensembles and:traces (2025)."
$ws.Range("A3").Value = "A2. Low annual flows within the:10-year periods (2025)."
$ws.Range("A4").Value = "B. From tree rings back to:1400 AD (2023)."
$ws.Range("A5").Value = "C. Collaborator choices in:immersive modeling:sessions (2021)."
$ws.Range("A6").Value = "D. 85%, 65%, and 50% of:2000 to 2018 average:flow (2022)."
$ws.Range("A7").Value = "E. Release to prevent:drawdown to 3,490 feet:(2024)."
$ws.Range("A8").Value = "F. Low Lake Powell:releases + gains through:Grand Canyon (2022)."

# Row 3 (A2 scenario) minimum/maximum flow values changed.
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 5

# Clear the lingering cell selection highlight left over from editing.
$ws.Range("A1").Select()
